# Insert a new weekly price record at row 243, shifting the existing
# rows 243:256 down to 244:257 (Excel's normal "insert row" behaviour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("243:243").Insert()

# Populate the newly inserted row with the new record's data. Columns
# A, B, C, E, F, G, H, I and R are identical for every row in this
# block, so copy them from the row directly below (old row 243, now
# shifted to row 244).
$ws.Range("A243").Value = $ws.Range("A244").Value()
$ws.Range("B243").Value = $ws.Range("B244").Value()
$ws.Range("C243").Value = $ws.Range("C244").Value()
$ws.Range("E243").Value = $ws.Range("E244").Value()
$ws.Range("F243").Value = $ws.Range("F244").Value()
$ws.Range("G243").Value = $ws.Range("G244").Value()
$ws.Range("H243").Value = $ws.Range("H244").Value()
$ws.Range("I243").Value = $ws.Range("I244").Value()
$ws.Range("R243").Value = $ws.Range("R244").Value()

# New record's own values.
$ws.Range("D243").Value = 44931
$ws.Range("J243").Value = 25
$ws.Range("K243").Value = 10000
$ws.Range("L243").Value = 10000
$ws.Range("M243").Value = 10000
$ws.Range("N243").Value = "$/docena de atados (12 kilos)"
$ws.Range("O243").Value = "Región de La Araucanía"
$ws.Range("P243").Value = 833
$ws.Range("Q243").Value = 12
